# context_local_storage.pptx
#   - refresh the cached "last saved" date field on the slide master
#   - widen two "CustomShape" boxes on slide 2 (room for the longer "->" arrow)
#   - replace "=" by "->" in the eight result lines on slide 2
#     ("replace = by -> for more readability, suggested by fred")

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide master: cached date field text 21-Feb-19 -> 1/15/2025
# ---------------------------------------------------------------------------
$master = $p.Slides.Item(2).Master
$dateShape = $master.Shapes.Item(2)
$dateShape.TextFrame.TextRange.Text = "1/15/2025"

# ---------------------------------------------------------------------------
# 2) Slide 2: widen the two rounded-rectangle outlines
#    (Shape.Width is a single-precision COM property, so the literal below is
#    picked so that it round-trips to exactly the target EMU value)
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)

$shape1 = $slide2.Shapes.Item(1)      # CustomShape 1 : 1645560 -> 1661400 EMU
$shape1.Width = 130.81893920898438

$shape12 = $slide2.Shapes.Item(12)    # CustomShape 12 : 5751720 -> 5881392 EMU
$shape12.Width = 463.1017761230469

# ---------------------------------------------------------------------------
# 3) Slide 2, CustomShape 12: "() = N" -> "() -> N" (8 occurrences)
# ---------------------------------------------------------------------------
$tr = $shape12.TextFrame.TextRange

$replacements = @(
    @("() = 1",              "() -> 1"),
    @("()  = 1",             "()  -> 1"),
    @("() = 1`t`t   ",       "() -> 1`t`t   "),
    @("()  = 1",             "()  -> 1"),
    @("()  = 1",             "()  -> 1"),
    @("() = 2",              "() -> 2"),
    @("() = 2",              "() -> 2"),
    @("() = 3",              "() -> 3")
)

$searchStart = 0
foreach ($pair in $replacements) {
    $oldStr = $pair[0]
    $newStr = $pair[1]
    $full = $tr.Text
    $idx = $full.IndexOf($oldStr, $searchStart)
    if ($idx -lt 0) {
        continue
    }
    $sub = $tr.Characters($idx + 1, $oldStr.Length)
    $sub.Text = $newStr
    $searchStart = $idx + $newStr.Length
}
